$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update header H1 to new shared string "Antarctic mass change (Gigatonnes)"
$ws.Range("H1").Value = "Antarctic mass change (Gigatonnes)"

# Widen column H to fit new header text
$ws.Columns.Item(8).ColumnWidth = 16

# Update formulas: H3:H20 first (keeps joining existing shared group),
# then H2 individually (stays a standalone formula, matching original layout)
$ws.Range("H3:H20").Formula = "=AVERAGEIF(A:A, G3, C:C)"
$ws.Range("H2").Formula = "=AVERAGEIF(A:A, G2, C:C)"

# Add new row 21 data: G21 = 2021, H21 = AVERAGEIF formula
$ws.Range("G21").Value = 2021
$ws.Range("H21").Formula = "=AVERAGEIF(A:A, G21, C:C)"

# Update selection to reflect extended range
$ws.Range("G1:H21").Select()

$wb.Save()
